$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "72.312.20"
$ws.Range("E2").Value = "  -0.15%  "
$ws.Range("D3").Value = "2.650.08"
$ws.Range("E3").Value = "  +0.11%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "591.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "175.29"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.53%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  -0.72%  "
$ws.Range("B9").Value = "LidoStakedEther"
$ws.Range("C9").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D9").Value = "2.649.10"
$ws.Range("E9").Value = "  +0.21%  "
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.172"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.46%  "
$ws.Range("E12").Value = "  -0.47%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.96"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.77%  "
$ws.Range("D14").Value = "3.134.83"
$ws.Range("E14").Value = "  +0.23%  "
$ws.Range("E15").Value = "  -2.09%  "
$ws.Range("D16").Value = "72.250.05"
$ws.Range("E16").Value = "  -0.11%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.06"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.72%  "
$ws.Range("D18").Value = "2.655.46"
$ws.Range("E18").Value = "  +0.31%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.27"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.28%  "
$ws.Range("E20").Value = "  +0.90%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "372.16"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.97%  "
$ws.Range("E22").Value = "  -0.18%  "
$ws.Range("E23").Value = "  +0.55%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.36"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.90%  "
$ws.Range("E25").Value = "  -0.22%  "
$ws.Range("E26").Value = "  -2.09%  "
$ws.Range("E27").Value = "  -2.87%  "
$ws.Range("D28").Value = "2.780.61"
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("E29").Value = "  -0.15%  "
$ws.Range("D30").Value = "0.0₃0964"
$ws.Range("E30").Value = "  +0.68%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.05"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.87%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "503.56"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.04%  "
$ws.Range("E33").Value = "  -1.94%  "
$ws.Range("E34").Value = "  -0.54%  "
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "161.80"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.35%  "
$ws.Range("B37").Value = "EthereumClassic"
$ws.Range("C37").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "19.37"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.26%  "
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.115"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.76%  "
$ws.Range("E39").Value = "  -1.02%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.36"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.34%  "
$ws.Range("E41").Value = "  -0.09%  "
$ws.Range("E42").Value = "  -6.31%  "
$ws.Range("E43").Value = "  -2.45%  "
$ws.Range("E44").Value = "  -3.30%  "
$ws.Range("E45").Value = "  -1.49%  "
$ws.Range("E46").Value = "  -0.57%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "153.65"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.39%  "
$ws.Range("E48").Value = "  +1.80%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.68"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.88%  "
$ws.Range("E50").Value = "  -0.58%  "
